$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Macro_taxonomy")

# Update the D column proportions for rows 3 and 4 (Urban ADO/STDRE -> 0.2)
$ws.Range("D3").Value = 0.2
$ws.Range("D4").Value = 0.2

# Row 6 used to be Wood/Urban/W/LWAL/1. It now becomes a new
# Block/Stone/Brick Urban CR/LFINF 0.1 row, and everything below
# shifts down by one (a new row is inserted to make room for the
# extra taxonomy line while the rest of the table shifts down).
$ws.Rows("6").Insert()

$ws.Range("A6").Value = "Block/Stone/Brick"
$ws.Range("B6").Value = "Urban"
$ws.Range("C6").Value = "CR/LFINF"
$ws.Range("D6").Value = 0.1

# Old row6 (Wood/Urban/W/LWAL/1) is now row7 - keep as-is.
$ws.Range("A7").Value = "Wood"
$ws.Range("B7").Value = "Urban"
$ws.Range("C7").Value = "W/LWAL"
$ws.Range("D7").Value = 1

# Old row7 (Block/Stone/Brick/Rural/MUR+CB/LWAL/0.2) is now row8,
# update its macro_taxonomy and proportion.
$ws.Range("C8").Value = "MUR+CB/LWAL"
$ws.Range("D8").Value = 0.2

# Old row8 (Block/Stone/Brick/Rural/MUR+ADO/LWAL/0.3) is now row9.
$ws.Range("C9").Value = "MUR+ADO/LWAL"
$ws.Range("D9").Value = 0.3

# Old row9 (Block/Stone/Brick/Rural/MUR+STRUB/LWAL/0.3) is now row10.
$ws.Range("C10").Value = "MUR+STRUB/LWAL"
$ws.Range("D10").Value = 0.3

# Old row10 (Block/Stone/Brick/Rural/MUR+CL/LWAL/0.2) is now row11.
$ws.Range("A11").Value = "Block/Stone/Brick"
$ws.Range("B11").Value = "Rural"
$ws.Range("C11").Value = "MUR+CL/LWAL"
$ws.Range("D11").Value = 0.2

# Old row11 (Wood/Rural/W/LWAL/1) is now row12.
$ws.Range("A12").Value = "Wood"
$ws.Range("B12").Value = "Rural"
$ws.Range("C12").Value = "W/LWAL"
$ws.Range("D12").Value = 1
